# Generate Report for Handback
#
# The localization-status report is refreshed once a handback completes:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    on the Overview sheet (zh-cn / de-de columns) and on each language
#    sheet's own Status column.
#  - The "Latest Handback DateTime" for zh-cn and de-de is stamped with the
#    new handback time.
#  - The "Error Detail" column is cleared now that the handback is in sync
#    (no more stale-handback-version error).
#  - The Status / Error Detail columns re-fit to their new (longer/shorter)
#    contents.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: E2 = zh-cn status, F2 = de-de status ---
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("K2").Value = "2016-09-02 15:01:58"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = 29.17
$wsZhCn.Columns.Item(16).ColumnWidth = 13.04

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("K2").Value = "2016-09-02 15:02:21"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = 29.17
$wsDeDe.Columns.Item(16).ColumnWidth = 13.04
